$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (T5YIFR) - latest date + present value shift
$ws.Range("N29").Value = "'2025-10-28"
$ws.Range("Q29").Value = 2.2

# Row 30 (T10YIE) - latest date + values shift right by one column
$ws.Range("N30").Value = "'2025-10-28"
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.29
$ws.Range("T30").Value = 2.3
$ws.Range("U30").Value = 2.29

# Row 37 (Case Shiller HPI, M/M % Delta) - latest date + recomputed values
$ws.Range("N37").Value = "'2025-08-01"
$ws.Range("Q37").Value = -0.002963719362060191
$ws.Range("R37").Value = -0.002056167749521443
$ws.Range("S37").Value = 0.0006848529828815675
$ws.Range("T37").Value = 0.004670871764716811
$ws.Range("U37").Value = 0.006620960680036703

# Row 38 (Case Shiller HPI, Y/Y % Delta) - latest date + recomputed values
$ws.Range("N38").Value = "'2025-08-01"
$ws.Range("Q38").Value = 0.01509936545139687
$ws.Range("R38").Value = 0.01637547210366313
$ws.Range("S38").Value = 0.019477605517784
$ws.Range("T38").Value = 0.02366620855788388
$ws.Range("U38").Value = 0.02826572083977447

# Row 48 (2y UST, DGS2) - latest date + values shift right by one column
$ws.Range("N48").Value = "'2025-10-27"
$ws.Range("S48").Value = 3.48
$ws.Range("U48").Value = 3.45

# Row 49 (5y UST, DGS5) - latest date + values shift right by one column
$ws.Range("N49").Value = "'2025-10-27"
$ws.Range("S49").Value = 3.61
$ws.Range("U49").Value = 3.56

# Row 50 (10y UST, DGS10) - latest date + values shift right by one column
$ws.Range("N50").Value = "'2025-10-27"
$ws.Range("Q50").Value = 4.01
$ws.Range("R50").Value = 4.02
$ws.Range("S50").Value = 4.01
$ws.Range("T50").Value = 3.97
$ws.Range("U50").Value = 3.98

# Row 52 (BAA, DBAA) - latest date + values shift right by one column
$ws.Range("N52").Value = "'2025-10-27"
$ws.Range("Q52").Value = 5.64
$ws.Range("S52").Value = 5.67
$ws.Range("T52").Value = 5.66
$ws.Range("U52").Value = 5.65
